# COVID19-Deaths.xlsx update: add three new timestamp columns (AU, AV, AW)
# for 2/13/20 10:00, 2/13/20 21:15 and 2/14/20 11:23, populate the death
# counts for every existing location row, and append a new row for
# "San Antonio, TX" (US).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cells for the three new report timestamps ------------
$ws.Range("AU1").Value = "2/13/20 10:00"
$ws.Range("AV1").Value = "2/13/20 21:15"
$ws.Range("AW1").Value = "2/14/20 11:23"

# Copy the header formatting (bold font, borders, centered alignment) from
# the preceding header cell (AT1) onto the new header cells.
$ws.Range("AT1").Copy()
$ws.Range("AU1:AW1").PasteSpecial(-4122)

# --- 2. Death counts for the three new timestamps, rows 2-74 ------------
$newCounts = @(
    @(5,6,6),
    @(3,3,3),
    @(3,4,5),
    @(0,0,0),
    @(2,2,2),
    @(2,2,2),
    @(2,2,2),
    @(1,1,1),
    @(4,4,4),
    @(3,3,3),
    @(9,11,11),
    @(10,11,11),
    @(1310,1426,1318),
    @(2,2,2),
    @(0,0,0),
    @(0,0,0),
    @(1,1,1),
    @(1,1,1),
    @(1,1,1),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(2,2,2),
    @(1,1,1),
    @(0,0,0),
    @(1,1,1),
    @(3,3,3),
    @(0,0,0),
    @(1,1,1),
    @(0,0,0),
    @(0,0,0),
    @(1,1,1),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(1,1,1),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(1,1,1),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0)
)

for ($i = 0; $i -lt $newCounts.Count; $i++) {
    $r = $i + 2
    $vals = $newCounts[$i]
    $ws.Cells.Item($r, 47).Value = $vals[0]
    $ws.Cells.Item($r, 48).Value = $vals[1]
    $ws.Cells.Item($r, 49).Value = $vals[2]
}

# --- 3. New row 75: San Antonio, TX (US) ---------------------------------
$ws.Range("A75").Value = "San Antonio, TX"
$ws.Range("B75").Value = "US"
$ws.Range("C75").Value = 29.4241
$ws.Range("D75").Value = -98.4936

for ($c = 5; $c -le 49; $c++) {
    $ws.Cells.Item(75, $c).Value = 0
}
